# "Accidental and Seismic combinations"
#
# 1) Actions sheet: factor/name tweaks (Accidental factors -> 1, Wind factors
#    tidy-up, Accidental/Seismic headers renamed, Seismic columns split into
#    "Seismic Horizontal" / "Seismic Vertical").
# 2) AccidentalLoadcases sheet: now envelopes Explosion + Impact (was
#    TLO Traffic / Wind / Thermal).
# 3) New SeismicLoadcases sheet (copied from AccidentalLoadcases) envelopes
#    Seismic Horizontal + Seismic Vertical.

$wb = $excel.ActiveWorkbook
$actions = $wb.Worksheets.Item("Actions")

# ---------------------------------------------------------------------
# 1) Actions sheet - numeric-only tweaks (no new text introduced)
# ---------------------------------------------------------------------

# Permanent Actions block - beneficial factor 0.95 -> 1
$actions.Range("B3").Value = 1
$actions.Range("B4").Value = 1
$actions.Range("B5").Value = 1

# Variable Actions block - Wind design factor tidy-up
$actions.Range("G5").Value = 1.5
$actions.Range("G6").Value = 1.5

# ---------------------------------------------------------------------
# Renames, in the same order the source file introduces them, so that
# the shared-string table is rebuilt in the same sequence as authored.
# ---------------------------------------------------------------------

# Seismic Actions block - split into Horizontal / Vertical
$actions.Range("S4").Value = "Seismic Vertical"
$actions.Range("S3").Value = "Seismic Horizontal"
$actions.Range("T3").Value = 1
$actions.Range("T4").Value = 1

# Accidental Actions block
$actions.Range("P4").Value = "Impact"
$actions.Range("Q3").Value = 1
$actions.Range("Q4").Value = 1

# ---------------------------------------------------------------------
# 2) AccidentalLoadcases sheet - now envelopes Explosion + Impact
# ---------------------------------------------------------------------
$acc = $wb.Worksheets.Item("AccidentalLoadcases")

$acc.Range("A2").Formula = "=Actions!P3"
$acc.Range("B2").Value = "Explosion Envelope"

$acc.Range("A3").Formula = "=Actions!P4"
$acc.Range("B3").Value = "Impact Envelope"

$acc.Range("A4").ClearContents()
$acc.Range("B4").ClearContents()

# ---------------------------------------------------------------------
# 3) New SeismicLoadcases sheet - clone AccidentalLoadcases's layout
# ---------------------------------------------------------------------
$acc.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$seismic = $wb.Worksheets.Item($wb.Worksheets.Count)
$seismic.Name = "SeismicLoadcases"

$seismic.Range("A2").Formula = "=Actions!S3"
$seismic.Range("B2").Value = "Seismic Envelope"

$seismic.Range("A3").Formula = "=Actions!S4"
$seismic.Range("B3").ClearContents()

$seismic.Columns.Item(1).ColumnWidth = 16

# Finally rename the Accidental/Seismic "Factor" headers - both Q2 and T2
# used the shared "Factor" text, so it only drops out of the table once
# both are replaced.
$actions.Range("Q2").Value = "AccFactor"
$actions.Range("T2").Value = "SeismicFactor"

# Cosmetic column widths for the newly-used columns around the Seismic block
$actions.Columns.Item(18).ColumnWidth = 4.83
$actions.Columns.Item(20).ColumnWidth = 13.17

# ---------------------------------------------------------------------
# Selections (mirrors the authored file's cursor positions)
# ---------------------------------------------------------------------
$perm = $wb.Worksheets.Item("PermanentLoadcases")
$perm.Activate()
$perm.Range("B2:C2").Select()

$var = $wb.Worksheets.Item("VariableLoadcases")
$var.Activate()
$var.Range("L12").Select()

$acc.Activate()
$acc.Range("C8").Select()

$seismic.Activate()
$seismic.Range("G12").Select()

# Actions must be the last-activated / selected sheet so it ends up the
# active tab on save, matching the source file.
$actions.Activate()
$actions.Range("S19").Select()
